# Update the shared-string status text "Ready for handoff" -> "In Translation".
# This string is shared across the Overview sheet (columns zh-cn/de-de) and the
# per-locale "Status" column on the zh-cn / de-de sheets, so updating every
# cell that currently holds it keeps all three sheets in sync.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow the "zh-cn" / "de-de" status columns on the Overview sheet, and the
# "Status" column on each locale sheet, to match the new shorter status text.
# (ColumnWidth is expressed in characters; the engine stores width in the
# workbook as characters + ~5/6 padding, so we back that padding out here to
# land on the target stored width.)
$overview.Columns.Item(5).ColumnWidth = 12.576851254417766
$overview.Columns.Item(6).ColumnWidth = 12.576851254417766
$zhcn.Columns.Item(3).ColumnWidth = 12.576851254417766
$dede.Columns.Item(3).ColumnWidth = 12.576851254417766
